# Updated all figures to use the same vaccination database.
# Update the age-group labels and dose percentages in the vaccination table
# so this sheet matches the shared vaccination database used elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to be treated as literal text so values such as "3%"
    # are not reinterpreted by Excel as a percentage number. Resetting the
    # style back to Normal afterwards avoids leaving a stray number format
    # applied to the cell itself.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 2: "0-2" age group
Set-TextValue "A2" "0-2 years old"
Set-TextValue "B2" "3%"
Set-TextValue "C2" "1%"

# Row 3: "3-11" age group
Set-TextValue "A3" "3-11 years old"
Set-TextValue "D3" "14%"

# Row 4: "12-17" age group
Set-TextValue "A4" "12-17 years old"
Set-TextValue "D4" "36%"
